$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 254, shifting the existing
# rows 254-330 down to 255-331 (dimension grows from A1:T330 to A1:T331).
$ws.Rows.Item(254).Insert()

# Populate the newly-inserted row 254 with the new Frutilla price record.
$ws.Cells.Item(254, 1).Value = 7
$ws.Cells.Item(254, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(254, 3).Value = "Ñuble"
$ws.Cells.Item(254, 4).Value = 44809
$ws.Cells.Item(254, 5).Value = 16
$ws.Cells.Item(254, 6).Value = "Fruta"
$ws.Cells.Item(254, 7).Value = 100101
$ws.Cells.Item(254, 8).Value = "Berries"
$ws.Cells.Item(254, 9).Value = 100112025
$ws.Cells.Item(254, 10).Value = "Frutilla"
$ws.Cells.Item(254, 11).Value = "Sin especificar"
$ws.Cells.Item(254, 12).Value = "Segunda"
$ws.Cells.Item(254, 13).Value = 30
$ws.Cells.Item(254, 14).Value = 15000
$ws.Cells.Item(254, 15).Value = 15000
$ws.Cells.Item(254, 16).Value = 15000
$ws.Cells.Item(254, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(254, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(254, 19).Value = 2143
$ws.Cells.Item(254, 20).Value = 7
